$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New settlement-ledger rows 32..43, appended after existing row 31 ---
# Each row: B=date label, C=status, D..I = raw ledger figures, J..M = formulas
# mirroring the pattern already used by rows 11..31.
#
# Row 31 carries the "blank" (no explicit style / style 34 / style 35)
# formatting that every new row should inherit, so we copy *only its
# formats* down onto each new row before filling in values/formulas --
# that keeps B..K free of any explicit style and L/M on styles 34/35,
# exactly like the existing unstyled rows 26..31.

$ws.Range("B31:M31").Copy()
for ($r = 32; $r -le 43; $r++) {
    $ws.Range("B$r`:M$r").PasteSpecial(-4122)
}

$rows = @(
    @{ r = 32; b = "2018.09.21 20:19:14"; d = 213883;  e = 368.999236;  f = 374333.661525625;  g = 355.991847325;  h = 588216.661525625;   i = 724.991083325 },
    @{ r = 33; b = "2018.09.21 20:19:14"; d = 213883;  e = 368.999236;  f = 374333.661525625;  g = 355.991847325;  h = 588216.661525625;   i = 724.991083325 },
    @{ r = 34; b = "2018.09.22 00:37:06"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 35; b = "2018.09.22 00:48:20"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 36; b = "2018.09.22 00:55:02"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 37; b = "2018.09.22 01:58:35"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 38; b = "2018.09.22 01:59:37"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 39; b = "2018.09.22 02:02:25"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 40; b = "2018.09.22 02:03:49"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 41; b = "2018.09.22 02:04:51"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 42; b = "2018.09.22 02:11:10"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 },
    @{ r = 43; b = "2018.09.22 02:12:18"; d = 344379;  e = 152.999236;  f = 251829.180347625;  g = 571.991631325;  h = 596208.1803476249;  i = 724.990867325 }
)

foreach ($row in $rows) {
    $r = $row.r
    $prev = $r - 1

    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = "initiation"
    $ws.Range("D$r").Value = $row.d
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Value = $row.f
    $ws.Range("G$r").Value = $row.g
    $ws.Range("H$r").Value = $row.h
    $ws.Range("I$r").Value = $row.i

    $ws.Range("J$r").Formula = "=IF(C$r=`"settlement`", H$r-H$prev, `"`")"
    $ws.Range("K$r").Formula = "=IF(C$r=`"settlement`", I$r-I$prev, `"`")"
    $ws.Range("L$r").Formula = "=IF(C$r=`"settlement`", J$r/H$prev, `"`")"
    $ws.Range("M$r").Formula = "=IF(C$r=`"settlement`", SUM(`$J`$11:J$r)/`$H`$11, `"`")"
}

# --- Fix up the three aggregate formulas at the top of the sheet so they
# cover the full (now 43-row) range instead of the old 31-row range. ---
$ws.Range("C5").Formula = "=SUM(J11:J43)"
$ws.Range("C6").Formula = "=SUM(K11:K43)"
$ws.Range("C7").Formula = "=M43"

Write-Output "edit complete"
